$d = $word.ActiveDocument

# The first two paragraphs of the document currently read:
#   Para 1: "On Pilgrimage" (italic) + ", " + <br/> + " " + "February =================="
#   Para 2: "By Dorothy Day" (bold)
#
# They need to become a pandoc-style title block:
#   Para 1 (styled "Title"): "February"
#   Para 2: "% Dorothy Day"

$p1 = $d.Paragraphs(1).Range
$p2 = $d.Paragraphs(2).Range
$targetRange = $d.Range($p1.Start, $p2.End)

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
              '<w:r><w:t xml:space="preserve">February</w:t></w:r>' +
            '</w:p>' +
            '<w:p>' +
              '<w:r><w:t xml:space="preserve">% Dorothy Day</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$targetRange.InsertXML($newXml)
